{"js": "// Apply the Coywolf Cash review text edits.\n// Title + SEO title (same text, appears twice): replaced with a new title.\n// \"What we like\" / \"What we don't like\" bullet lists: replaced with new copy.\n// Meta description (italic paragraph): replaced with a new description.\n\nconst replacements = [\n  [\n    \"Play Coywolf Cash Free: High-Volatility Slot Review\",\n    \"Play Coywolf Cash Free & Win Big | Online Slot Game\",\n  ],\n  [\n    \"High volatility with up to 5000x payout\",\n    \"Wide betting range for both casual and daring players\",\n  ],\n  [\n    \"Attractive graphics and design elements\",\n    \"High volatility with chances to win big\",\n  ],\n  [\n    \"Multiple bonus features for more chances to win\",\n    \"Stunning graphics and immersive design elements\",\n  ],\n  [\n    \"Compatible across various devices\",\n    \"Several bonus features for more opportunities to win\",\n  ],\n  [\n    \"May not appeal to players looking for low volatility\",\n    \"Limited number of paylines\",\n  ],\n  [\n    \"May not offer enough variety in gameplay\",\n    \"May not appeal to players who prefer low volatility slots\",\n  ],\n  [\n    \"Experience the American wilderness with Coywolf Cash. Read our review and play for free to win big payouts on this high-volatility online slot game.\",\n    \"Review of Coywolf Cash, an online slot game with stunning graphics and bonus features. Play for free and win big online.\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Coywolf Cash review text edits via Find/Replace.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"Play Coywolf Cash Free: High-Volatility Slot Review\", \"Play Coywolf Cash Free & Win Big | Online Slot Game\"),\n    @(\"High volatility with up to 5000x payout\", \"Wide betting range for both casual and daring players\"),\n    @(\"Attractive graphics and design elements\", \"High volatility with chances to win big\"),\n    @(\"Multiple bonus features for more chances to win\", \"Stunning graphics and immersive design elements\"),\n    @(\"Compatible across various devices\", \"Several bonus features for more opportunities to win\"),\n    @(\"May not appeal to players looking for low volatility\", \"Limited number of paylines\"),\n    @(\"May not offer enough variety in gameplay\", \"May not appeal to players who prefer low volatility slots\"),\n    @(\"Experience the American wilderness with Coywolf Cash. Read our review and play for free to win big payouts on this high-volatility online slot game.\", \"Review of Coywolf Cash, an online slot game with stunning graphics and bonus features. Play for free and win big online.\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $find\n    $rng.Find.Replacement.Text = $replace\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n"}
